# Add a "Save" column (H) to the s_vals sheet, mirroring the header style
# used by the existing header row and filling row 2 with 0.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell H1 = "Save", matching the header formatting of G1 ("sum").
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# New data cell H2 = 0 (plain number, no special style - like the rest of row 2).
$ws.Range("H2").Value = 0
